$wb = $excel.ActiveWorkbook
$snap = $wb.Worksheets.Item("snapshot")

# Insert two new rows at position 31 (shifts old rows 31-39 down to 33-41)
$snap.Range("A31:A32").EntireRow.Insert()

# New row 31: SPR / Korostelev Nikita
# New row 32: SYuL / Alalykin Danil
$snap.Cells.Item(31,1).Value = 'СПР'
$snap.Cells.Item(31,2).Value = 'Спартак'
$snap.Cells.Item(31,3).Value = 'spartak'
$snap.Cells.Item(31,4).Value = 'Коростелёв Никита'
$snap.Cells.Item(31,5).Value = "'35"
$snap.Cells.Item(31,6).Value = 'нападающий'
$snap.Cells.Item(31,7).Value = "'22149"
$snap.Cells.Item(31,8).Value = '1369_СПР_коростелевникита'
$snap.Cells.Item(31,9).Value = 'injured_active'
$snap.Cells.Item(31,10).Value = 'https://www.khl.ru/clubs/spartak/team/'
$snap.Cells.Item(31,11).Value = '2025-12-11T03:01:56.972285+00:00'

$snap.Cells.Item(32,1).Value = 'СЮЛ'
$snap.Cells.Item(32,2).Value = 'Салават Юлаев'
$snap.Cells.Item(32,3).Value = 'salavat_yulaev'
$snap.Cells.Item(32,4).Value = 'Алалыкин Данил'
$snap.Cells.Item(32,5).Value = "'61"
$snap.Cells.Item(32,6).Value = 'нападающий'
$snap.Cells.Item(32,7).Value = "'34493"
$snap.Cells.Item(32,8).Value = '1369_СЮЛ_алалыкинданил'
$snap.Cells.Item(32,9).Value = 'injured_active'
$snap.Cells.Item(32,10).Value = 'https://www.khl.ru/clubs/salavat_yulaev/team/'
$snap.Cells.Item(32,11).Value = '2025-12-11T03:02:10.084372+00:00'

# Refresh scraped_at (column K) for every row -- this run re-scraped all teams
$snap.Cells.Item(2,11).Value = '2025-12-11T03:01:18.469656+00:00'
$snap.Cells.Item(3,11).Value = '2025-12-11T03:01:18.469695+00:00'
$snap.Cells.Item(4,11).Value = '2025-12-11T03:01:18.469716+00:00'
$snap.Cells.Item(5,11).Value = '2025-12-11T03:01:20.786893+00:00'
$snap.Cells.Item(6,11).Value = '2025-12-11T03:01:20.786922+00:00'
$snap.Cells.Item(7,11).Value = '2025-12-11T03:01:20.786940+00:00'
$snap.Cells.Item(8,11).Value = '2025-12-11T03:01:23.575330+00:00'
$snap.Cells.Item(9,11).Value = '2025-12-11T03:01:25.874413+00:00'
$snap.Cells.Item(10,11).Value = '2025-12-11T03:01:28.211919+00:00'
$snap.Cells.Item(11,11).Value = '2025-12-11T03:01:30.932641+00:00'
$snap.Cells.Item(12,11).Value = '2025-12-11T03:01:35.977679+00:00'
$snap.Cells.Item(13,11).Value = '2025-12-11T03:01:35.977709+00:00'
$snap.Cells.Item(14,11).Value = '2025-12-11T03:01:38.227288+00:00'
$snap.Cells.Item(15,11).Value = '2025-12-11T03:01:41.000370+00:00'
$snap.Cells.Item(16,11).Value = '2025-12-11T03:01:43.707714+00:00'
$snap.Cells.Item(17,11).Value = '2025-12-11T03:01:46.523675+00:00'
$snap.Cells.Item(18,11).Value = '2025-12-11T03:01:46.523703+00:00'
$snap.Cells.Item(19,11).Value = '2025-12-11T03:01:48.793953+00:00'
$snap.Cells.Item(20,11).Value = '2025-12-11T03:01:48.793982+00:00'
$snap.Cells.Item(21,11).Value = '2025-12-11T03:01:48.793999+00:00'
$snap.Cells.Item(22,11).Value = '2025-12-11T03:01:51.473256+00:00'
$snap.Cells.Item(23,11).Value = '2025-12-11T03:01:51.473285+00:00'
$snap.Cells.Item(24,11).Value = '2025-12-11T03:01:51.473302+00:00'
$snap.Cells.Item(25,11).Value = '2025-12-11T03:01:51.473320+00:00'
$snap.Cells.Item(26,11).Value = '2025-12-11T03:01:51.473336+00:00'
$snap.Cells.Item(27,11).Value = '2025-12-11T03:01:54.258740+00:00'
$snap.Cells.Item(28,11).Value = '2025-12-11T03:01:54.258988+00:00'
$snap.Cells.Item(29,11).Value = '2025-12-11T03:01:54.259024+00:00'
$snap.Cells.Item(30,11).Value = '2025-12-11T03:01:56.972255+00:00'
$snap.Cells.Item(33,11).Value = '2025-12-11T03:02:12.831779+00:00'
$snap.Cells.Item(34,11).Value = '2025-12-11T03:02:12.831812+00:00'
$snap.Cells.Item(35,11).Value = '2025-12-11T03:02:12.831832+00:00'
$snap.Cells.Item(36,11).Value = '2025-12-11T03:02:15.140860+00:00'
$snap.Cells.Item(37,11).Value = '2025-12-11T03:02:15.140888+00:00'
$snap.Cells.Item(38,11).Value = '2025-12-11T03:02:17.966492+00:00'
$snap.Cells.Item(39,11).Value = '2025-12-11T03:02:17.966521+00:00'
$snap.Cells.Item(40,11).Value = '2025-12-11T03:02:20.661942+00:00'
$snap.Cells.Item(41,11).Value = '2025-12-11T03:02:20.661970+00:00'

$newinj = $wb.Worksheets.Item("new_injured")

$newinj.Cells.Item(2,1).Value = 'СПР'
$newinj.Cells.Item(2,2).Value = 'Спартак'
$newinj.Cells.Item(2,3).Value = 'Коростелёв Никита'
$newinj.Cells.Item(2,4).Value = '1369_СПР_коростелевникита'
$newinj.Cells.Item(2,5).Value = 'INJURED_NEW'
$newinj.Cells.Item(2,6).Value = '2025-12-11T11:02:21.170929+08:00'
$newinj.Cells.Item(2,7).Value = "'2025-12-11"

$newinj.Cells.Item(3,1).Value = 'СЮЛ'
$newinj.Cells.Item(3,2).Value = 'Салават Юлаев'
$newinj.Cells.Item(3,3).Value = 'Алалыкин Данил'
$newinj.Cells.Item(3,4).Value = '1369_СЮЛ_алалыкинданил'
$newinj.Cells.Item(3,5).Value = 'INJURED_NEW'
$newinj.Cells.Item(3,6).Value = '2025-12-11T11:02:21.170929+08:00'
$newinj.Cells.Item(3,7).Value = "'2025-12-11"
